$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "65.747.59"
Set-TextValue "E2" "  +0.49%  "

# Row 3
Set-TextValue "D3" "2.672.98"
Set-TextValue "E3" "  +0.97%  "

# Row 4
Set-TextValue "E4" "  -0.03%  "

# Row 5
Set-TextValue "D5" "601.17"
Set-TextValue "E5" "  -0.55%  "

# Row 6
Set-TextValue "D6" "157.11"
Set-TextValue "E6" "  +0.76%  "

# Row 7
Set-TextValue "E7" "  -0.01%  "

# Row 8
Set-TextValue "D8" "0.618"
Set-TextValue "E8" "  +5.37%  "

# Row 9
Set-TextValue "E9" "  +1.01%  "

# Row 10
Set-TextValue "D10" "5.94"
Set-TextValue "E10" "  +1.47%  "

# Row 11
Set-TextValue "E11" "  -0.18%  "

# Row 12
Set-TextValue "E12" "  -0.24%  "

# Row 13
Set-TextValue "D13" "29.50"
Set-TextValue "E13" "  -0.89%  "

# Row 14
Set-TextValue "D14" "0.0000197"
Set-TextValue "E14" "  +1.73%  "

# Row 15
Set-TextValue "D15" "3.150.17"
Set-TextValue "E15" "  +0.81%  "

# Row 16
Set-TextValue "D16" "65.560.62"
Set-TextValue "E16" "  +0.53%  "

# Row 17
Set-TextValue "D17" "2.668.04"
Set-TextValue "E17" "  +0.92%  "

# Row 18
Set-TextValue "D18" "12.64"
Set-TextValue "E18" "  -0.50%  "

# Row 19
Set-TextValue "D19" "4.82"
Set-TextValue "E19" "  -1.08%  "

# Row 20
Set-TextValue "E20" "  +1.42%  "

# Row 21
Set-TextValue "D21" "352.37"
Set-TextValue "E21" "  -1.62%  "

# Row 22
Set-TextValue "E22" "  -0.09%  "

# Row 23
Set-TextValue "D23" "69.87"
Set-TextValue "E23" "  +0.33%  "

# Row 24
Set-TextValue "E24" "  +6.14%  "

# Row 25
Set-TextValue "D25" "9.80"
Set-TextValue "E25" "  +4.30%  "

# Row 26
Set-TextValue "D26" "1.62"
Set-TextValue "E26" "  -4.62%  "

# Row 27
Set-TextValue "E27" "  +1.72%  "

# Row 28
Set-TextValue "E28" "  -1.45%  "

# Row 29
Set-TextValue "D29" "8.12"
Set-TextValue "E29" "  +0.30%  "

# Row 30
Set-TextValue "D30" "543.63"
Set-TextValue "E30" "  +3.76%  "

# Row 31
Set-TextValue "D31" "0.993"
Set-TextValue "E31" "  -0.61%  "

# Row 32
Set-TextValue "D32" "2.17"
Set-TextValue "E32" "  +0.22%  "

# Row 33
Set-TextValue "E33" "  +0.11%  "

# Row 34
Set-TextValue "D34" "6.58"
Set-TextValue "E34" "  +4.28%  "

# Row 35
Set-TextValue "D35" "5.47"
Set-TextValue "E35" "  -0.63%  "

# Row 36
Set-TextValue "E36" "  -1.49%  "

# Row 37
Set-TextValue "D37" "20.43"
Set-TextValue "E37" "  -0.90%  "

# Row 38
Set-TextValue "D38" "0.999"
Set-TextValue "E38" "  -0.03%  "

# Row 39
Set-TextValue "D39" "158.33"
Set-TextValue "E39" "  -2.50%  "

# Row 40
Set-TextValue "E40" "  -0.79%  "

# Row 41
Set-TextValue "E41" "  +0.04%  "

# Row 42
Set-TextValue "D42" "42.74"

# Row 43
Set-TextValue "D43" "164.94"
Set-TextValue "E43" "  +0.00%  "

# Row 44
Set-TextValue "E44" "  -0.60%  "

# Row 45
Set-TextValue "D45" "0.0614"
Set-TextValue "E45" "  +1.24%  "

# Row 46
Set-TextValue "E46" "  -0.57%  "

# Row 47
Set-TextValue "D47" "23.24"
Set-TextValue "E47" "  +1.64%  "

# Row 48
Set-TextValue "E48" "  -0.91%  "

# Row 49
Set-TextValue "E49" "  -1.21%  "

# Row 50
Set-TextValue "E50" "  +3.25%  "

# Row 51
Set-TextValue "D51" "20.28"
Set-TextValue "E51" "  +3.14%  "
